# Commit message: "Removed spaces from .csv file"
# Observed effect: a new "Width" column is populated in column C (which
# previously held an unused "Brightness" header with empty data cells).
# Each row's width is determined by its Species (column B) value - this
# mirrors a per-species line-width lookup that was merged in from the
# source .csv after its column headers/whitespace were cleaned up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Species -> Width lookup (line width value per species)
$widthBySpecies = @{
    "O I"   = 0.25
    "Cl I"  = 0.26
    "S III" = 0.4
    "S II"  = 0.31
    "S I"   = 0.23
}

# Re-purpose column C: header "Brightness" -> "Width"
$ws.Range("C1").Value = "Width"

# Find the last used data row (column A holds Wavelength values down to row 34)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $species = $ws.Cells.Item($r, 2).Value2
    if ($widthBySpecies.ContainsKey($species)) {
        $cell = $ws.Cells.Item($r, 3)
        $cell.Value = $widthBySpecies[$species]
        $cell.Style = "Normal"
    }
}

# Update the active selection to match the author's final cursor position
[void]$ws.Range("C34").Select()
